# Update res_line pl_mw results for the 380 kV case (Case_1_90).
# For each data row (rows 2-25), refresh columns B, D:I and L:N with the
# newly computed line active-power-loss results; columns A, C, J, K, O are
# unchanged (index / zero columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8322655858124506
$row2DI = New-Object 'object[,]' 1,6
$row2DI[0,0] = 0.00466880969363892
$row2DI[0,1] = 0.5789671946912733
$row2DI[0,2] = 0.9090700550331974
$row2DI[0,3] = 0.8196547599940516
$row2DI[0,4] = 0.7153065298483909
$row2DI[0,5] = 0.7897106584390654
$ws.Range("D2:I2").Value = $row2DI
$row2LMN = New-Object 'object[,]' 1,3
$row2LMN[0,0] = 0.2820181567121693
$row2LMN[0,1] = 0.2429056420267983
$row2LMN[0,2] = 1.534823580371238
$ws.Range("L2:N2").Value = $row2LMN

$ws.Range("B3").Value = 0.783107803365084
$row3DI = New-Object 'object[,]' 1,6
$row3DI[0,0] = 0.004660741958552705
$row3DI[0,1] = 0.5445172120720798
$row3DI[0,2] = 0.8547575287720264
$row3DI[0,3] = 0.7577624574712445
$row3DI[0,4] = 0.6923911283652444
$row3DI[0,5] = 0.8033681632161898
$ws.Range("D3:I3").Value = $row3DI
$row3LMN = New-Object 'object[,]' 1,3
$row3LMN[0,0] = 0.257329694569421
$row3LMN[0,1] = 0.225739417751349
$row3LMN[0,2] = 1.535672379390149
$ws.Range("L3:N3").Value = $row3LMN

$ws.Range("B4").Value = 0.7531503350731725
$row4DI = New-Object 'object[,]' 1,6
$row4DI[0,0] = 0.004663326268536494
$row4DI[0,1] = 0.5232130336073624
$row4DI[0,2] = 0.8221097482335011
$row4DI[0,3] = 0.7203873233417255
$row4DI[0,4] = 0.6788620256523643
$row4DI[0,5] = 0.8122311290833917
$ws.Range("D4:I4").Value = $row4DI
$row4LMN = New-Object 'object[,]' 1,3
$row4LMN[0,0] = 0.2422605457202138
$row4LMN[0,1] = 0.2152721607178805
$row4LMN[0,2] = 1.536788252829467
$ws.Range("L4:N4").Value = $row4LMN

$ws.Range("B5").Value = 0.7409998991072939
$row5DI = New-Object 'object[,]' 1,6
$row5DI[0,0] = 0.004666244512929296
$row5DI[0,1] = 0.5144931407772759
$row5DI[0,2] = 0.8089797074441236
$row5DI[0,3] = 0.7053120244728746
$row5DI[0,4] = 0.6734839374936996
$row5DI[0,5] = 0.8159627345177549
$ws.Range("D5:I5").Value = $row5DI
$row5LMN = New-Object 'object[,]' 1,3
$row5LMN[0,0] = 0.2361422577095738
$row5LMN[0,1] = 0.2110251354267731
$row5LMN[0,2] = 1.537393068067516
$ws.Range("L5:N5").Value = $row5LMN

$ws.Range("B6").Value = 0.7389858218823235
$row6DI = New-Object 'object[,]' 1,6
$row6DI[0,0] = 0.004666840850314102
$row6DI[0,1] = 0.5130428949890131
$row6DI[0,2] = 0.806809939610929
$row6DI[0,3] = 0.7028181000219149
$row6DI[0,4] = 0.6725990449261303
$row6DI[0,5] = 0.8165896013583485
$ws.Range("D6:I6").Value = $row6DI
$row6LMN = New-Object 'object[,]' 1,3
$row6LMN[0,0] = 0.2351276794387331
$row6LMN[0,1] = 0.2103210397090436
$row6LMN[0,2] = 1.53750257823539
$ws.Range("L6:N6").Value = $row6LMN

$ws.Range("B7").Value = 0.7529862363427071
$row7DI = New-Object 'object[,]' 1,6
$row7DI[0,0] = 0.004663358113329963
$row7DI[0,1] = 0.5230955891507136
$row7DI[0,2] = 0.8219319691251314
$row7DI[0,3] = 0.7201833865125877
$row7DI[0,4] = 0.6787889489579015
$row7DI[0,5] = 0.8122809696150544
$ws.Range("D7:I7").Value = $row7DI
$row7LMN = New-Object 'object[,]' 1,3
$row7LMN[0,0] = 0.2421779412406551
$row7LMN[0,1] = 0.2152148088891437
$row7LMN[0,2] = 1.536795801158391
$ws.Range("L7:N7").Value = $row7LMN

$ws.Range("B8").Value = 0.8152695859491246
$row8DI = New-Object 'object[,]' 1,6
$row8DI[0,0] = 0.004664445081456137
$row8DI[0,1] = 0.5671202595685401
$row8DI[0,2] = 0.8901965716210896
$row8DI[0,3] = 0.7981827604934324
$row8DI[0,4] = 0.7072924406344896
$row8DI[0,5] = 0.7943206509173732
$ws.Range("D8:I8").Value = $row8DI
$row8LMN = New-Object 'object[,]' 1,3
$row8LMN[0,0] = 0.2734869345614328
$row8LMN[0,1] = 0.2369716819273222
$row8LMN[0,2] = 1.534993135342262
$ws.Range("L8:N8").Value = $row8LMN

$ws.Range("B9").Value = 0.939171220586644
$row9DI = New-Object 'object[,]' 1,6
$row9DI[0,0] = 0.004727769982135754
$row9DI[0,1] = 0.6522575607046974
$row9DI[0,2] = 1.029714146113022
$row9DI[0,3] = 0.9562260061970278
$row9DI[0,4] = 0.7675261763798176
$row9DI[0,5] = 0.7628918429592026
$ws.Range("D9:I9").Value = $row9DI
$row9LMN = New-Object 'object[,]' 1,3
$row9LMN[0,0] = 0.3356003118632316
$row9LMN[0,1] = 0.2802112171940792
$row9LMN[0,2] = 1.536154608198061
$ws.Range("L9:N9").Value = $row9LMN

$ws.Range("B10").Value = 1.031252393707121
$row10DI = New-Object 'object[,]' 1,6
$row10DI[0,0] = 0.004813535364398547
$row10DI[0,1] = 0.714098226033272
$row10DI[0,2] = 1.135804447183034
$row10DI[0,3] = 1.075608514509952
$row10DI[0,4] = 0.814493605582129
$row10DI[0,5] = 0.742118425336395
$ws.Range("D10:I10").Value = $row10DI
$row10LMN = New-Object 'object[,]' 1,3
$row10LMN[0,0] = 0.3816839008875093
$row10LMN[0,1] = 0.312327830058436
$row10LMN[0,2] = 1.539843567417293
$ws.Range("L10:N10").Value = $row10LMN

$ws.Range("B11").Value = 1.073366075168735
$row11DI = New-Object 'object[,]' 1,6
$row11DI[0,0] = 0.004861496786407571
$row11DI[0,1] = 0.7420815639944749
$row11DI[0,2] = 1.184878013289932
$row11DI[0,3] = 1.130665385245834
$row11DI[0,4] = 0.8364649156988264
$row11DI[0,5] = 0.7331725372843154
$ws.Range("D11:I11").Value = $row11DI
$row11LMN = New-Object 'object[,]' 1,3
$row11LMN[0,0] = 0.402748958035346
$row11LMN[0,1] = 0.3270140998094249
$row11LMN[0,2] = 1.542131968117701
$ws.Range("L11:N11").Value = $row11LMN

$ws.Range("B12").Value = 1.089345262696213
$row12DI = New-Object 'object[,]' 1,6
$row12DI[0,0] = 0.004880979987280654
$row12DI[0,1] = 0.7526570495825951
$row12DI[0,2] = 1.20358009027737
$row12DI[0,3] = 1.151624499239858
$row12DI[0,4] = 0.8448731117883312
$row12DI[0,5] = 0.7298575782917007
$ws.Range("D12:I12").Value = $row12DI
$row12LMN = New-Object 'object[,]' 1,3
$row12LMN[0,0] = 0.4107404674194584
$row12LMN[0,1] = 0.3325862918456366
$row12LMN[0,2] = 1.543085802645749
$ws.Range("L12:N12").Value = $row12LMN

$ws.Range("B13").Value = 1.085902463378602
$row13DI = New-Object 'object[,]' 1,6
$row13DI[0,0] = 0.004876724609626848
$row13DI[0,1] = 0.7503803701751934
$row13DI[0,2] = 1.199546937291387
$row13DI[0,3] = 1.147105633895649
$row13DI[0,4] = 0.8430583194244434
$row13DI[0,5] = 0.7305682799032978
$ws.Range("D13:I13").Value = $row13DI
$row13LMN = New-Object 'object[,]' 1,3
$row13LMN[0,0] = 0.4090187014246567
$row13LMN[0,1] = 0.3313857416469901
$row13LMN[0,2] = 1.542876503804166
$ws.Range("L13:N13").Value = $row13LMN

$ws.Range("B14").Value = 1.074680062643154
$row14DI = New-Object 'object[,]' 1,6
$row14DI[0,0] = 0.004863073000390017
$row14DI[0,1] = 0.7429520396206613
$row14DI[0,2] = 1.186414247500096
$row14DI[0,3] = 1.13238747949697
$row14DI[0,4] = 0.8371548905508064
$row14DI[0,5] = 0.7328983565048155
$ws.Range("D14:I14").Value = $row14DI
$row14LMN = New-Object 'object[,]' 1,3
$row14LMN[0,0] = 0.4034061308350658
$row14LMN[0,1] = 0.3274723113717073
$row14LMN[0,2] = 1.542208694476244
$ws.Range("L14:N14").Value = $row14LMN

$ws.Range("B15").Value = 1.067810117055501
$row15DI = New-Object 'object[,]' 1,6
$row15DI[0,0] = 0.004854884100872425
$row15DI[0,1] = 0.7383992160962123
$row15DI[0,2] = 1.178385654706375
$row15DI[0,3] = 1.123386621517682
$row15DI[0,4] = 0.8335503795782984
$row15DI[0,5] = 0.7343350626670571
$ws.Range("D15:I15").Value = $row15DI
$row15LMN = New-Object 'object[,]' 1,3
$row15LMN[0,0] = 0.3999701759758523
$row15LMN[0,1] = 0.3250766271787455
$row15LMN[0,2] = 1.541810992308669
$ws.Range("L15:N15").Value = $row15LMN

$ws.Range("B16").Value = 1.028504636667151
$row16DI = New-Object 'object[,]' 1,6
$row16DI[0,0] = 0.004810583958253289
$row16DI[0,1] = 0.7122664757293364
$row16DI[0,2] = 1.132613907146634
$row16DI[0,3] = 1.072025710450049
$row16DI[0,4] = 0.8130700046678214
$row16DI[0,5] = 0.7427132176295759
$ws.Range("D16:I16").Value = $row16DI
$row16LMN = New-Object 'object[,]' 1,3
$row16LMN[0,0] = 0.3803092943152819
$row16LMN[0,1] = 0.3113695698407355
$row16LMN[0,2] = 1.539706253256895
$ws.Range("L16:N16").Value = $row16LMN

$ws.Range("B17").Value = 1.004449198251336
$row17DI = New-Object 'object[,]' 1,6
$row17DI[0,0] = 0.004785723229641548
$row17DI[0,1] = 0.6961969924619495
$row17DI[0,2] = 1.104744051326009
$row17DI[0,3] = 1.040711135903848
$row17DI[0,4] = 0.8006618104291476
$row17DI[0,5] = 0.7479821651971745
$ws.Range("D17:I17").Value = $row17DI
$row17LMN = New-Object 'object[,]' 1,3
$row17LMN[0,0] = 0.3682739842500098
$row17LMN[0,1] = 0.302980161960555
$row17LMN[0,2] = 1.538571018624168
$ws.Range("L17:N17").Value = $row17LMN

$ws.Range("B18").Value = 0.99063441683478
$row18DI = New-Object 'object[,]' 1,6
$row18DI[0,0] = 0.004772263958983558
$row18DI[0,1] = 0.6869402847043062
$row18DI[0,2] = 1.088790399745136
$row18DI[0,3] = 1.022770182606678
$row18DI[0,4] = 0.793581882515781
$row18DI[0,5] = 0.7510601570011222
$ws.Range("D18:I18").Value = $row18DI
$row18LMN = New-Object 'object[,]' 1,3
$row18LMN[0,0] = 0.3613611334388622
$row18LMN[0,1] = 0.2981619831015507
$row18LMN[0,2] = 1.537975532370808
$ws.Range("L18:N18").Value = $row18LMN

$ws.Range("B19").Value = 0.9859606461784551
$row19DI = New-Object 'object[,]' 1,6
$row19DI[0,0] = 0.004767850093760018
$row19DI[0,1] = 0.6838037201289069
$row19DI[0,2] = 1.083401816272925
$row19DI[0,3] = 1.016707691746234
$row19DI[0,4] = 0.7911944876799737
$row19DI[0,5] = 0.7521104534436649
$ws.Range("D19:I19").Value = $row19DI
$row19LMN = New-Object 'object[,]' 1,3
$row19LMN[0,0] = 0.3590221982594812
$row19LMN[0,1] = 0.2965318693281915
$row19LMN[0,2] = 1.537783796666233
$ws.Range("L19:N19").Value = $row19LMN

$ws.Range("B20").Value = 1.007007744955558
$row20DI = New-Object 'object[,]' 1,6
$row20DI[0,0] = 0.004788282512631525
$row20DI[0,1] = 0.6979090609260936
$row20DI[0,2] = 1.107702927467329
$row20DI[0,3] = 1.04403732531722
$row20DI[0,4] = 0.8019767835075413
$row20DI[0,5] = 0.7474163663955196
$ws.Range("D20:I20").Value = $row20DI
$row20LMN = New-Object 'object[,]' 1,3
$row20LMN[0,0] = 0.3695541757395802
$row20LMN[0,1] = 0.3038724860234154
$row20LMN[0,2] = 1.538685921869742
$ws.Range("L20:N20").Value = $row20LMN

$ws.Range("B21").Value = 1.07797550094125
$row21DI = New-Object 'object[,]' 1,6
$row21DI[0,0] = 0.004867046671428454
$row21DI[0,1] = 0.7451344945985028
$row21DI[0,2] = 1.190268390522249
$row21DI[0,3] = 1.136707545075865
$row21DI[0,4] = 0.8388864702052672
$row21DI[0,5] = 0.7322119833057696
$ws.Range("D21:I21").Value = $row21DI
$row21LMN = New-Object 'object[,]' 1,3
$row21LMN[0,0] = 0.4050542814241567
$row21LMN[0,1] = 0.3286214881941163
$row21LMN[0,2] = 1.542402481727379
$ws.Range("L21:N21").Value = $row21LMN

$ws.Range("B22").Value = 1.124541159041257
$row22DI = New-Object 'object[,]' 1,6
$row22DI[0,0] = 0.004926241327463998
$row22DI[0,1] = 0.7758755577492025
$row22DI[0,2] = 1.244924224364951
$row22DI[0,3] = 1.197916838867741
$row22DI[0,4] = 0.8635233178192152
$row22DI[0,5] = 0.7226986087058105
$ws.Range("D22:I22").Value = $row22DI
$row22LMN = New-Object 'object[,]' 1,3
$row22LMN[0,0] = 0.4283409497535899
$row22LMN[0,1] = 0.3448593954928469
$row22LMN[0,2] = 1.545339905221823
$ws.Range("L22:N22").Value = $row22LMN

$ws.Range("B23").Value = 1.099671597876466
$row23DI = New-Object 'object[,]' 1,6
$row23DI[0,0] = 0.004893930243767386
$row23DI[0,1] = 0.7594797285815531
$row23DI[0,2] = 1.215689114887255
$row23DI[0,3] = 1.165188518577764
$row23DI[0,4] = 0.8503267641458478
$row23DI[0,5] = 0.7277372613881399
$ws.Range("D23:I23").Value = $row23DI
$row23LMN = New-Object 'object[,]' 1,3
$row23LMN[0,0] = 0.4159045948300957
$row23LMN[0,1] = 0.3361872015006142
$row23LMN[0,2] = 1.543725784620861
$ws.Range("L23:N23").Value = $row23LMN

$ws.Range("B24").Value = 1.005850979176671
$row24DI = New-Object 'object[,]' 1,6
$row24DI[0,0] = 0.0047871228673948
$row24DI[0,1] = 0.6971350913090362
$row24DI[0,2] = 1.106365004494734
$row24DI[0,3] = 1.042533361496453
$row24DI[0,4] = 0.8013821169689663
$row24DI[0,5] = 0.7476720120343927
$ws.Range("D24:I24").Value = $row24DI
$row24LMN = New-Object 'object[,]' 1,3
$row24LMN[0,0] = 0.3689753812072354
$row24LMN[0,1] = 0.3034690507186468
$row24LMN[0,2] = 1.538633795986854
$ws.Range("L24:N24").Value = $row24LMN

$ws.Range("B25").Value = 0.9054662701643679
$row25DI = New-Object 'object[,]' 1,6
$row25DI[0,0] = 0.004703889204547451
$row25DI[0,1] = 0.6293516832676858
$row25DI[0,2] = 0.9913514011765017
$row25DI[0,3] = 0.9129093644720285
$row25DI[0,4] = 0.7507598203254418
$row25DI[0,5] = 0.7709876630632353
$ws.Range("D25:I25").Value = $row25DI
$row25LMN = New-Object 'object[,]' 1,3
$row25LMN[0,0] = 0.318719055786147
$row25LMN[0,1] = 0.2684524634247083
$row25LMN[0,2] = 1.535340482427529
$ws.Range("L25:N25").Value = $row25LMN

